$d = $word.ActiveDocument

# Locate the paragraph that contains the stale M2Doc version-mismatch
# warning (inserted by M2Doc validation when the template/runtime
# versions differ) and strip it back down to the single empty run that
# was there before validation added the warning runs.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*M2Doc version mismatch*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End
    # Keep the paragraph mark (last character); delete everything else
    # in the paragraph so only the original empty run remains.
    $r = $d.Range($pStart, $pEnd - 1)
    $r.Delete()
}
